# Updated symbol list on Fri Feb 10 23:46:39 UTC 2023 with GitHub Actions
# Applies refreshed Price (column D) and Volume(1h) (column E) values
# for the crypto rows that changed in this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPrice; E = newVolume } (only columns that changed are included)
$updates = [ordered]@{
    2  = @{ D = "306.11";        E = "-0.25%" }
    3  = @{ D = "40.86";         E = "3.73%"  }
    4  = @{ D = "5.106";         E = "2.29%"  }
    5  = @{ D = "0.07609";       E = "-1.44%" }
    6  = @{ D = "4.264";         E = "-0.02%" }
    7  = @{ D = "1.616";         E = "1.10%"  }
    8  = @{             E = "-5.28%" }
    9  = @{ D = "0.9066";        E = "-0.84%" }
    10 = @{ D = "0.1015";        E = "0.52%"  }
    11 = @{ D = "0.1750";        E = "1.54%"  }
    12 = @{ D = "0.09068";       E = "0.69%"  }
    13 = @{ D = "0.04263";       E = "-4.04%" }
    14 = @{ D = "0.1055";        E = "-0.30%" }
    15 = @{ D = "0.001226";      E = "-3.78%" }
    16 = @{ D = "0.005848";      E = "3.55%"  }
    17 = @{             E = "-0.30%" }
    19 = @{ D = "6.541";         E = "-7.30%" }
    20 = @{ D = "0.1355";        E = "-0.61%" }
    21 = @{             E = "-4.73%" }
    22 = @{ D = "0.04177";       E = "0.93%"  }
    23 = @{ D = "0.001226";      E = "2.28%"  }
    24 = @{ D = "0.004074";      E = "-0.10%" }
    25 = @{ D = "0.0001301";     E = "6.29%"  }
    26 = @{ D = "0.0003005";     E = "0.54%"  }
    38 = @{ D = "0.02378";       E = "1.34%"  }
    39 = @{ D = "0.05127";       E = "0.07%"  }
    40 = @{ D = "0.007773";      E = "-2.67%" }
    41 = @{ D = "0.1297";        E = "-2.27%" }
    42 = @{ D = "0.007044";      E = "-3.73%" }
    43 = @{ D = "0.001919";      E = "-4.73%" }
    44 = @{ D = "0.008437" }
    45 = @{ D = "0.3327";        E = "0.53%"  }
    46 = @{ D = "0.00006351";    E = "-5.15%" }
    47 = @{ D = "0.00000000750"; E = "-0.40%" }
    48 = @{ D = "0.007109";      E = "109.61%" }
    49 = @{ D = "0.004398";      E = "6.88%"  }
    50 = @{ D = "0.00002099";    E = "-0.40%" }
    51 = @{ D = "0.0001999";     E = "-0.40%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        $origStyle = $cell.Style
        # Force text storage so values like "306.11" / "-0.25%" stay exact strings
        # (matching the inline-string cells already used in this sheet) instead of
        # being auto-converted by Excel into numbers / percentages.
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
        # Restore the cell's original style/number format so formatting is unchanged.
        $cell.Style = $origStyle
    }
}
